# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" detail sheet (fund holdings) right after
#    "2021-Q4" (i.e. right before "总计").
# 2. Re-create the "总计" sheet after it so it keeps its trailing position
#    and picks up the next sheetId (the "总计" -> sheetId 7 / "2022-Q1" ->
#    sheetId 6 split matches the real edit).
# 3. Populate "2022-Q1" with the per-fund holdings table.
# 4. Populate "总计" with the previous summary rows plus a new leading row
#    for 2022-Q1.
#
# Formatting is carried over by Range.Copy()-ing cells that already have the
# workbook's "header / index column" style (bold, centered, thin border)
# instead of re-declaring fonts/borders by hand, so no new cellXf entries are
# introduced and plain data cells are left with the default (unstyled) xf —
# exactly like every other quarter sheet in this workbook.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$lastDetailSheet = $wb.Worksheets.Item("2021-Q4")
$blankFormatCell = $lastDetailSheet.Cells.Item(1, 26)   # untouched cell -> default/no style

# Delete the old "总计" sheet FIRST so the two sheets (re)created below are
# assigned sheetId 6 then 7, matching the target state.
$wb.Worksheets.Item("总计").Delete() | Out-Null

# ---------------------------------------------------------------------------
# Step 1: create "2022-Q1" at the end (after "2021-Q4").
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$q1.Name = "2022-Q1"

# Header row (B1:H1) and index column (A2:A16) pick up the bold/centered/
# thin-border style already used for those roles on every other sheet.
$lastDetailSheet.Range("B1").Copy($q1.Range("B1:H1")) | Out-Null
$lastDetailSheet.Range("A2").Copy($q1.Range("A2:A16")) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: re-create "总计" after "2022-Q1".
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$total.Name = "总计"

$lastDetailSheet.Range("B1").Copy($total.Range("B1:D1")) | Out-Null
$lastDetailSheet.Range("A2").Copy($total.Range("A2:A7")) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: fill in the "2022-Q1" fund-holdings table.
# ---------------------------------------------------------------------------
$fundRows = @(
    @("161834", "银华鑫锐灵活配置混合（LOF）",         "67.33", "81.90", "2.05", "1.3803", 9),
    @("501022", "银华鑫盛灵活配置混合（LOF）",         "61.98", "79.75", "2.01", "1.2458", 8),
    @("002666", "前海开源沪港深创新成长灵活配置混合A", "11.96", "81.64", "5.84", "0.6985", 9),
    @("260112", "景顺长城能源基建混合",                 "16.49", "60.89", "2.51", "0.4139", 6),
    @("009782", "富国兴泉回报12个月持有期混合A",       "6.18",  "70.06", "4.56", "0.2818", 4),
    @("012370", "银华鑫利一年持有期混合型证券投资基金", "11.01", "80.06", "2.02", "0.2224", 8),
    @("002667", "前海开源沪港深创新成长灵活配置混合C", "3.25",  "81.64", "5.84", "0.1898", 9),
    @("011046", "富国优质企业混合A",                     "8.18",  "71.23", "2.27", "0.1857", 9),
    @("001306", "中欧永裕混合A",                         "4.48",  "86.33", "3.29", "0.1474", 9),
    @("009783", "富国兴泉回报12个月持有期混合C",       "2.45",  "70.06", "4.56", "0.1117", 4),
    @("013678", "富国信享回报12个月持有期混合A",       "9.49",  "27.59", "1.04", "0.0987", 10),
    @("005732", "富国臻选成长灵活配置混合",             "2.45",  "64.81", "3.98", "0.0975", 3),
    @("013679", "富国信享回报12个月持有期混合C",       "2.76",  "27.59", "1.04", "0.0287", 10),
    @("001307", "中欧永裕混合C",                         "0.35",  "86.33", "3.29", "0.0115", 9),
    @("011047", "富国优质企业混合C",                     "0.48",  "71.23", "2.27", "0.0109", 9)
)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B-G hold numeric-looking text (fund codes need leading zeros kept,
# and the scale/weight figures are stored as text in the source data), so
# force text format before writing the values...
$q1.Range("B2:G16").NumberFormat = "@"

$rowIdx = 2
$idxValue = 0
foreach ($fund in $fundRows) {
    $q1.Cells.Item($rowIdx, 1).Value = $idxValue
    $q1.Cells.Item($rowIdx, 2).Value = $fund[0]
    $q1.Cells.Item($rowIdx, 3).Value = $fund[1]
    $q1.Cells.Item($rowIdx, 4).Value = $fund[2]
    $q1.Cells.Item($rowIdx, 5).Value = $fund[3]
    $q1.Cells.Item($rowIdx, 6).Value = $fund[4]
    $q1.Cells.Item($rowIdx, 7).Value = $fund[5]
    $q1.Cells.Item($rowIdx, 8).Value = $fund[6]
    $rowIdx = $rowIdx + 1
    $idxValue = $idxValue + 1
}

# ...then strip the now-redundant explicit "@" number-format/style back off
# those cells (formats only; values are untouched) so they end up as plain
# unstyled text cells, matching every other quarter sheet.
$blankFormatCell.Copy() | Out-Null
$q1.Range("B2:G16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# ---------------------------------------------------------------------------
# Step 4: rebuild the "总计" summary sheet with the new leading 2022-Q1 row.
# ---------------------------------------------------------------------------
$summaryRows = @(
    @("2022-Q1", 15, 5.12),
    @("2021-Q4", 12, 4.16),
    @("2021-Q3", 61, 23.94),
    @("2021-Q2", 23, 4.48),
    @("2021-Q1", 19, 3.05),
    @("2020-Q4", 11, 1.83)
)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$rowIdx = 2
$idxValue = 0
foreach ($entry in $summaryRows) {
    $total.Cells.Item($rowIdx, 1).Value = $idxValue
    $total.Cells.Item($rowIdx, 2).Value = $entry[0]
    $total.Cells.Item($rowIdx, 3).Value = $entry[1]
    $total.Cells.Item($rowIdx, 4).Value = $entry[2]
    $rowIdx = $rowIdx + 1
    $idxValue = $idxValue + 1
}

$wb.Worksheets.Item("2021-Q4").Activate()
